$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric values in column B
$ws.Range("B2").Value = 37.39
$ws.Range("B3").Value = 69.49
$ws.Range("B5").Value = 0.627
$ws.Range("B6").Value = 0.327
$ws.Range("B7").Value = 0.327
$ws.Range("B8").Value = 0.627

# Move the active cell / selection to D12 (cosmetic, matches authored diff)
$ws.Range("D12").Select()
